# Update the cryptos price/volume table (rows 2-51) on the active sheet to
# reflect the latest coinranking.com snapshot. A couple of same-rank ties
# swapped places (Litecoin/Dai, Kaspa/FirstDigitalUSD, dogwifhat/USDe), so
# those rows also get new Coin name + Link values.
#
# Price strings like "595.10" would be auto-parsed into numbers by Excel's
# Range.Value setter (losing the trailing zero / turning "1.00" into 1), so
# those are written with a leading apostrophe to force text, matching the
# original inline-string cell content exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '66.988.51'
$ws.Range("E2").Value = '  -2.17%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '2.665.26'
$ws.Range("E3").Value = '  -0.92%  '

# Row 4: TetherUSD
$ws.Range("E4").Value = '  +0.02%  '

# Row 5: BNB
$ws.Range("D5").Value = '''595.10'
$ws.Range("E5").Value = '  -0.52%  '

# Row 6: Solana
$ws.Range("D6").Value = '''163.60'
$ws.Range("E6").Value = '  +2.73%  '

# Row 7: USDC
$ws.Range("E7").Value = '  +0.00%  '

# Row 8: XRP
$ws.Range("E8").Value = '  +0.46%  '

# Row 9: LidoStakedEther
$ws.Range("D9").Value = '2.666.29'
$ws.Range("E9").Value = '  -0.85%  '

# Row 10: Dogecoin
$ws.Range("E10").Value = '  +0.97%  '

# Row 11: TRON
$ws.Range("D11").Value = '''0.158'
$ws.Range("E11").Value = '  +0.66%  '

# Row 12: Cardano
$ws.Range("E12").Value = '  -0.60%  '

# Row 13: Toncoin
$ws.Range("E13").Value = '  -1.89%  '

# Row 14: Avalanche
$ws.Range("D14").Value = '''27.69'
$ws.Range("E14").Value = '  -2.00%  '

# Row 15: WrappedliquidstakedEther2.0
$ws.Range("D15").Value = '3.164.82'

# Row 16: ShibaInu
$ws.Range("D16").Value = '''0.0000182'
$ws.Range("E16").Value = '  -2.64%  '

# Row 17: WrappedBTC
$ws.Range("D17").Value = '66.966.36'
$ws.Range("E17").Value = '  -2.15%  '

# Row 18: WrappedEther
$ws.Range("D18").Value = '2.680.01'
$ws.Range("E18").Value = '  -1.70%  '

# Row 19: Chainlink
$ws.Range("D19").Value = '''11.59'
$ws.Range("E19").Value = '  -1.97%  '

# Row 20: BitcoinCash
$ws.Range("D20").Value = '''360.74'
$ws.Range("E20").Value = '  -1.27%  '

# Row 21: Uniswap
$ws.Range("D21").Value = '''7.49'
$ws.Range("E21").Value = '  -1.42%  '

# Row 22: Polkadot
$ws.Range("E22").Value = '  -3.38%  '

# Row 23: NEARProtocol
$ws.Range("D23").Value = '''4.79'
$ws.Range("E23").Value = '  -1.91%  '

# Row 24: SuiNetwork
$ws.Range("D24").Value = '''2.02'
$ws.Range("E24").Value = '  -4.66%  '

# Row 25: Dai
$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D25").Value = '''1.00'
$ws.Range("E25").Value = '  +0.17%  '

# Row 26: Litecoin
$ws.Range("B26").Value = 'Litecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D26").Value = '''71.17'
$ws.Range("E26").Value = '  -5.15%  '

# Row 27: Aptos
$ws.Range("D27").Value = '''10.04'
$ws.Range("E27").Value = '  -0.47%  '

# Row 28: WrappedeETH
$ws.Range("E28").Value = '  -1.41%  '

# Row 29: Binance-PegBSC-USD
$ws.Range("D29").Value = '''1.00'
$ws.Range("E29").Value = '  -0.13%  '

# Row 30: PEPE
$ws.Range("D30").Value = '''0.0000102'
$ws.Range("E30").Value = '  -2.10%  '

# Row 31: Bittensor
$ws.Range("D31").Value = '''553.08'
$ws.Range("E31").Value = '  -4.42%  '

# Row 32: InternetComputer(DFINITY)
$ws.Range("E32").Value = '  -3.84%  '

# Row 33: Fetch.AI
$ws.Range("E33").Value = '  -3.23%  '

# Row 34: PancakeSwap
$ws.Range("E34").Value = '  -0.72%  '

# Row 35: Kaspa
$ws.Range("B35").Value = 'Kaspa'
$ws.Range("C35").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D35").Value = '''0.128'
$ws.Range("E35").Value = '  -1.68%  '

# Row 36: FirstDigitalUSD
$ws.Range("B36").Value = 'FirstDigitalUSD'
$ws.Range("C36").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D36").Value = '''0.999'
$ws.Range("E36").Value = '  +0.00%  '

# Row 37: ImmutableX
$ws.Range("E37").Value = '  -4.72%  '

# Row 38: EthereumClassic
$ws.Range("D38").Value = '''19.40'
$ws.Range("E38").Value = '  -2.91%  '

# Row 39: Monero
$ws.Range("D39").Value = '''155.80'
$ws.Range("E39").Value = '  -3.59%  '

# Row 40: PolygonEcosystemToken
$ws.Range("D40").Value = '''0.372'
$ws.Range("E40").Value = '  -2.14%  '

# Row 41: RenderToken
$ws.Range("D41").Value = '''5.26'
$ws.Range("E41").Value = '  -2.35%  '

# Row 42: Stacks
$ws.Range("E42").Value = '  -4.46%  '

# Row 43: WhiteBITCoin
$ws.Range("E43").Value = '  +0.19%  '

# Row 44: dogwifhat
$ws.Range("B44").Value = 'dogwifhat'
$ws.Range("C44").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D44").Value = '''2.52'
$ws.Range("E44").Value = '  -4.35%  '

# Row 45: USDe
$ws.Range("B45").Value = 'USDe'
$ws.Range("C45").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D45").Value = '''1.00'
$ws.Range("E45").Value = '  +0.06%  '

# Row 46: OKB
$ws.Range("D46").Value = '''40.23'
$ws.Range("E46").Value = '  -0.65%  '

# Row 47: BabyDogeCoin
$ws.Range("D47").Value = '0.0₆0296'
$ws.Range("E47").Value = '  -6.31%  '

# Row 48: ARBITRUM
$ws.Range("D48").Value = '''0.584'
$ws.Range("E48").Value = '  -2.72%  '

# Row 49: Aave
$ws.Range("D49").Value = '''152.45'
$ws.Range("E49").Value = '  -3.77%  '

# Row 50: Filecoin
$ws.Range("D50").Value = '''3.82'
$ws.Range("E50").Value = '  -2.69%  '

# Row 51: Optimism
$ws.Range("E51").Value = '  -2.50%  '
